# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (keep as date serial, matching existing date formatting of A2)
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("A2").Value = $epoch.AddDays(45974)

# Hourly prices 0h-1h ... 23h-24h (columns B..Y)
$ws.Range("B2").Value = 24.52
$ws.Range("C2").Value = 21.56
$ws.Range("D2").Value = 17.26
$ws.Range("E2").Value = 14.23
$ws.Range("F2").Value = 13.9
$ws.Range("G2").Value = 17.01
$ws.Range("H2").Value = 30.19
$ws.Range("I2").Value = 45.87
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 43.03
$ws.Range("L2").Value = 30.39
$ws.Range("M2").Value = 21.7
$ws.Range("N2").Value = 22.26
$ws.Range("O2").Value = 20.99
$ws.Range("P2").Value = 21.81
$ws.Range("Q2").Value = 30.12
$ws.Range("R2").Value = 38.85
$ws.Range("S2").Value = 45.04
$ws.Range("T2").Value = 53.96
$ws.Range("U2").Value = 53.12
$ws.Range("V2").Value = 49.19
$ws.Range("W2").Value = 44.53
$ws.Range("X2").Value = 33.38
$ws.Range("Y2").Value = 22.71

# Daily average
$ws.Range("Z2").Value = 31.9

# 4h slot max window + price
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 47.74

# 2h slot windows + prices (first then second)
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 53.54
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 46.86

# Minimum price slot
$ws.Range("AG2").Value = "0h-23h"
